$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of (cell reference, new value) pairs reflecting the latest cryptos scrape.
# NumberFormat is forced to Text ("@") before assignment so numeric-looking strings
# (e.g. "253.51") are stored verbatim as text instead of being auto-converted to
# floating point numbers, matching the original text cells in the workbook.
$updates = @(
    ,("D2", "37.147.01")
    ,("E2", "  -0.41%  ")
    ,("D3", "2.081.06")
    ,("E3", "  -0.59%  ")
    ,("E4", "  +0.12%  ")
    ,("D5", "253.51")
    ,("E5", "  +1.19%  ")
    ,("D6", "0.674")
    ,("E6", "  +2.28%  ")
    ,("D7", "58.91")
    ,("E7", "  +13.91%  ")
    ,("D8", "0.999")
    ,("E8", "  -0.10%  ")
    ,("D9", "0.391")
    ,("E9", "  +5.00%  ")
    ,("D10", "61.74")
    ,("E10", "  +0.09%  ")
    ,("E11", "  +7.08%  ")
    ,("E12", "  +2.61%  ")
    ,("D13", "16.03")
    ,("E13", "  +6.69%  ")
    ,("D14", "2.380.21")
    ,("E14", "  -0.93%  ")
    ,("D15", "0.818")
    ,("E15", "  -1.85%  ")
    ,("E16", "  +8.20%  ")
    ,("D17", "2.069.76")
    ,("E17", "  -1.59%  ")
    ,("D18", "37.111.05")
    ,("E18", "  -0.44%  ")
    ,("D19", "15.59")
    ,("E19", "  +10.18%  ")
    ,("D20", "74.58")
    ,("E20", "  +3.23%  ")
    ,("D21", "0.0₃0925")
    ,("E21", "  +9.95%  ")
    ,("E22", "  +4.70%  ")
    ,("D23", "239.59")
    ,("E23", "  -0.26%  ")
    ,("E24", "  +0.07%  ")
    ,("D25", "2.42")
    ,("E25", "  -1.72%  ")
    ,("D26", "2.31")
    ,("E26", "  +14.90%  ")
    ,("D27", "169.95")
    ,("E27", "  -0.58%  ")
    ,("D28", "9.32")
    ,("E28", "  +0.96%  ")
    ,("D29", "20.34")
    ,("E29", "  -1.53%  ")
    ,("E30", "  +2.74%  ")
    ,("E31", "  +7.72%  ")
    ,("D32", "1.13")
    ,("E32", "  +5.90%  ")
    ,("D33", "0.0637")
    ,("E33", "  +4.25%  ")
    ,("D34", "4.45")
    ,("E34", "  +8.78%  ")
    ,("D35", "0.0919")
    ,("E35", "  +0.35%  ")
    ,("E36", "  +0.00%  ")
    ,("E37", "  -0.28%  ")
    ,("B38", "Cronos")
    ,("C38", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro")
    ,("D38", "0.118")
    ,("E38", "  +28.87%  ")
    ,("B39", "WEMIXToken")
    ,("C39", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix")
    ,("D39", "1.77")
    ,("E39", "  -3.91%  ")
    ,("D40", "1.37")
    ,("E40", "  +2.57%  ")
    ,("D41", "4.71")
    ,("E41", "  +30.93%  ")
    ,("B42", "InjectiveProtocol")
    ,("C42", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj")
    ,("D42", "17.92")
    ,("E42", "  -2.41%  ")
    ,("B43", "VeChain")
    ,("C43", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet")
    ,("D43", "0.0227")
    ,("E43", "  +1.56%  ")
    ,("D44", "1.17")
    ,("E44", "  +0.55%  ")
    ,("D45", "99.01")
    ,("E45", "  +0.13%  ")
    ,("E46", "  +2.15%  ")
    ,("D47", "4.47")
    ,("E47", "  +12.95%  ")
    ,("E48", "  +8.65%  ")
    ,("D49", "2.97")
    ,("E49", "  -0.63%  ")
    ,("D50", "1.306.35")
    ,("E50", "  -0.94%  ")
    ,("D51", "6.95")
    ,("E51", "  -0.58%  ")
)

foreach ($pair in $updates) {
    $rng = $ws.Range($pair[0])
    $rng.NumberFormat = "@"
    $rng.Value = $pair[1]
}
